$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the "Status" shared text wherever it appears: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "Ready for handoff"

# Latest HO Xliff Generate Date: "2016-08-31 08:29:34" -> "2016-08-31 08:32:27"
$wsOverview.Range("G2").Value = "2016-08-31 08:32:27"
$wsOverview.Range("G3").Value = "2016-08-31 08:32:27"

# Priority "ht" -> "mt" (shared across zh-cn and de-de rows)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# Latest Handoff Datetime "2016-08-31 08:29:21" -> "2016-08-31 08:32:14" (zh-cn only)
$wsZhCn.Range("H2").Value = "2016-08-31 08:32:14"
$wsZhCn.Range("H3").Value = "2016-08-31 08:32:14"

# de-de Latest Handoff Datetime shares the same string as Overview's "Latest HO Xliff Generate Date"
# ("2016-08-31 08:29:34" -> "2016-08-31 08:32:27"); set explicitly to keep it in sync.
$wsDeDe.Range("H2").Value = "2016-08-31 08:32:27"
$wsDeDe.Range("H3").Value = "2016-08-31 08:32:27"

# Error Detail column (P) for row 2 in both zh-cn and de-de sheets
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c65fc68726df8c4b9cf0daf0113bfc1fbf00aa9/e2e/2b6ee2ab-b090-4b82-a03a-4b2d2f0d859d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15f0a9d1bd608121b3c94007fc2d5ad89a880829/e2e/2b6ee2ab-b090-4b82-a03a-4b2d2f0d859d.md."
$wsZhCn.Range("P2").Value = $errorDetail
$wsDeDe.Range("P2").Value = $errorDetail

# Column width changes (narrower date columns, wider Error Detail column)
$wsOverview.Range("E1").ColumnWidth = 16.333333333333332
$wsOverview.Range("F1").ColumnWidth = 16.333333333333332

$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("P1").ColumnWidth = 39.166666666666664

$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("P1").ColumnWidth = 39.166666666666664
